$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timesheet: append three new logged entries below the existing data
# (rows 2-28), mirroring the date/hours/comment layout already used.

# Row 29: 2020-07-20, 1h, "XSLConstructor w osobnym programie"
$ws.Cells.Item(29, 1).Value = 44032
$ws.Cells.Item(29, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(29, 2).Value = 1
$ws.Cells.Item(29, 3).Value = "XSLConstructor w osobnym programie"

# Row 30: 2020-07-20, 1h, "RabbitMQ"
$ws.Cells.Item(30, 1).Value = 44032
$ws.Cells.Item(30, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(30, 2).Value = 1
$ws.Cells.Item(30, 3).Value = "RabbitMQ"

# Row 31: 2020-07-22, 4h, "RabbitMQ"
$ws.Cells.Item(31, 1).Value = 44034
$ws.Cells.Item(31, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(31, 2).Value = 4
$ws.Cells.Item(31, 3).Value = "RabbitMQ"

# E2 is =SUM(B:B); recalculates automatically to include the new hours.

# Move the selection to the next empty row under column C, same as the
# author leaving the cursor ready for the next entry.
$ws.Range("C32").Select()
